# Scheduled runner update: refresh Leve profit calcs (currentAveragePrice /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H-N)
# across the ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 830.5
$ws.Range("I4").Value = 746
$ws.Range("K4").Value = 746
$ws.Range("M4").Value = -632
$ws.Range("H18").Value = 2509.5
$ws.Range("I18").Value = 1420
$ws.Range("K18").Value = 1420
$ws.Range("M18").Value = -1136
$ws.Range("H19").Value = 968.8333
$ws.Range("I19").Value = 1487.75
$ws.Range("J19").Value = 709.375
$ws.Range("K19").Value = 1487.75
$ws.Range("L19").Value = 709.375
$ws.Range("M19").Value = -1312.75
$ws.Range("N19").Value = -1059.375
$ws.Range("H58").Value = 1655.5555
$ws.Range("I58").Value = 1292.1
$ws.Range("K58").Value = 3876.3
$ws.Range("M58").Value = -3726.3
$ws.Range("H112").Value = 1989.2667
$ws.Range("I112").Value = 2000
$ws.Range("J112").Value = 1986.5834
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 5959.7502
$ws.Range("M112").Value = -4892
$ws.Range("N112").Value = -8175.7502
$ws.Range("H116").Value = 3981.5
$ws.Range("I116").Value = 3981.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3981.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -539.5
$ws.Range("N116").ClearContents() | Out-Null
$ws.Range("H125").Value = 2942.5
$ws.Range("I125").Value = 3385
$ws.Range("J125").Value = 2500
$ws.Range("K125").Value = 30465
$ws.Range("L125").Value = 22500
$ws.Range("M125").Value = -28005
$ws.Range("N125").Value = -27420
$ws.Range("H127").Value = 626
$ws.Range("I127").Value = 626
$ws.Range("K127").Value = 1878
$ws.Range("M127").Value = 3082
$ws.Range("H132").Value = 9223.666999999999
$ws.Range("I132").Value = 10716.19
$ws.Range("K132").Value = 32148.57
$ws.Range("M132").Value = -29618.57

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents() | Out-Null
$ws.Range("H61").Value = 5435.7144
$ws.Range("I61").Value = 4110.8184
$ws.Range("K61").Value = 4110.8184
$ws.Range("M61").Value = -3898.8184
$ws.Range("H74").Value = 6532.3335
$ws.Range("I74").Value = 6549.125
$ws.Range("K74").Value = 6549.125
$ws.Range("M74").Value = -5675.125
$ws.Range("H77").Value = 6532.3335
$ws.Range("I77").Value = 6549.125
$ws.Range("K77").Value = 32745.625
$ws.Range("M77").Value = -28377.625
$ws.Range("H102").Value = 2910.2104
$ws.Range("I102").Value = 1349.5714
$ws.Range("J102").Value = 7280
$ws.Range("K102").Value = 1349.5714
$ws.Range("L102").Value = 7280
$ws.Range("M102").Value = 272.4286
$ws.Range("N102").Value = -10524
$ws.Range("H136").Value = 5435.7144
$ws.Range("I136").Value = 4110.8184
$ws.Range("K136").Value = 12332.4552
$ws.Range("M136").Value = -9782.4552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 20000.334
$ws.Range("J6").Value = 29000.5
$ws.Range("L6").Value = 29000.5
$ws.Range("N6").Value = -29226.5
$ws.Range("H16").Value = 1900
$ws.Range("I16").Value = 1900
$ws.Range("K16").Value = 1900
$ws.Range("M16").Value = -1613
$ws.Range("H19").Value = 478.33334
$ws.Range("I19").Value = 417.5
$ws.Range("J19").Value = 600
$ws.Range("K19").Value = 417.5
$ws.Range("L19").Value = 600
$ws.Range("M19").Value = -247.5
$ws.Range("N19").Value = -940
$ws.Range("H22").Value = 1978.2307
$ws.Range("I22").Value = 1644.2
$ws.Range("J22").Value = 2187
$ws.Range("K22").Value = 1644.2
$ws.Range("L22").Value = 2187
$ws.Range("M22").Value = -1294.2
$ws.Range("N22").Value = -2887
$ws.Range("H24").Value = 478.33334
$ws.Range("I24").Value = 417.5
$ws.Range("J24").Value = 600
$ws.Range("K24").Value = 417.5
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = -247.5
$ws.Range("N24").Value = -940
$ws.Range("H33").Value = 2730
$ws.Range("I33").Value = 2730
$ws.Range("K33").Value = 2730
$ws.Range("M33").Value = -2351
$ws.Range("H36").Value = 550
$ws.Range("I36").Value = 550
$ws.Range("K36").Value = 550
$ws.Range("M36").Value = -162
$ws.Range("H40").Value = 550
$ws.Range("I40").Value = 550
$ws.Range("K40").Value = 550
$ws.Range("M40").Value = -390
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3877
$ws.Range("N86").ClearContents() | Out-Null
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -19384
$ws.Range("H113").Value = 1900
$ws.Range("I113").Value = 1900
$ws.Range("K113").Value = 1900
$ws.Range("M113").Value = 270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 108.4375
$ws.Range("I12").Value = 33.25
$ws.Range("J12").Value = 133.5
$ws.Range("K12").Value = 99.75
$ws.Range("L12").Value = 400.5
$ws.Range("M12").Value = 73.25
$ws.Range("N12").Value = -746.5
$ws.Range("H117").Value = 1554
$ws.Range("J117").Value = 1554
$ws.Range("L117").Value = 4662
$ws.Range("N117").Value = -11546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 337.91306
$ws.Range("I2").Value = 192.57143
$ws.Range("K2").Value = 192.57143
$ws.Range("M2").Value = -79.57142999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2669.5908
$ws.Range("I82").Value = 1451.4166
$ws.Range("J82").Value = 4131.4
$ws.Range("K82").Value = 1451.4166
$ws.Range("L82").Value = 4131.4
$ws.Range("M82").Value = -1090.4166
$ws.Range("N82").Value = -4853.4
$ws.Range("H85").Value = 2669.5908
$ws.Range("I85").Value = 1451.4166
$ws.Range("J85").Value = 4131.4
$ws.Range("K85").Value = 1451.4166
$ws.Range("L85").Value = 4131.4
$ws.Range("M85").Value = -203.4166
$ws.Range("N85").Value = -6627.4
$ws.Range("H136").Value = 3873.75
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3873.75
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 11621.25
$ws.Range("M136").ClearContents() | Out-Null
$ws.Range("N136").Value = -16721.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 54249.5
$ws.Range("I27").Value = 53999
$ws.Range("J27").Value = 54333
$ws.Range("K27").Value = 53999
$ws.Range("L27").Value = 54333
$ws.Range("M27").Value = -53930
$ws.Range("N27").Value = -54471
$ws.Range("H113").Value = 896.63635
$ws.Range("I113").Value = 818.2222
$ws.Range("K113").Value = 2454.6666
$ws.Range("M113").Value = -284.6666
$ws.Range("H132").Value = 663.8889
$ws.Range("J132").Value = 1022.5
$ws.Range("L132").Value = 3067.5
$ws.Range("N132").Value = -8127.5
